# Added Case Auto Entitlement Test cases
# - Update the Service Contract creation/fetch sample code in F2/G2
#   (contract renamed from 'SCON RS_1022' to 'SCON BO Regression', and the
#   fetch query no longer filters on Name).
# - Add a new "CaseTriggerEnable" column (M) with its header and sample code.
# - Resize the new columns and refresh the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2 / G2: new contract name, simplified fetch query -------------------
$ws.Range("F2").Value = "SVMXC__Service_Contract__c service_contract = New SVMXC__Service_Contract__c (Name = 'SCON BO Regression', SVMXC__Active__c = true , SVMXC__All_Contacts_Covered__c = true , SVMXC__Company__c = '001q000000kxZfw');insert service_contract;"
$ws.Range("G2").Value = "Select Name , Id from SVMXC__Service_Contract__c where Createdby.Id = '005q0000003GGfP' Order by CreatedDate DESC Limit 1"

# --- M1 / M2: new "CaseTriggerEnable" column -------------------------------
$ws.Range("M1").Value = "CaseTriggerEnable"
$ws.Range("M2").Value = 'SVMXC.COMM_Utils_ManageSettings cums = new SVMXC.COMM_Utils_ManageSettings(); 
cums.SVMX_getInventoryProcessSteps(''{"propertyKey": "SVMXC_CASE_Trigger1","status":"Enabled","orgId":"00Dq0000000933B"}''); '


# --- Column widths for the two new columns ---------------------------------
$ws.Columns.Item(13).ColumnWidth = 58.25
$ws.Columns.Item(14).ColumnWidth = 31.25

# --- Refresh the view: scroll so column F is leftmost, select I2 ----------
$ws.Range("F1").Select()
$ws.Range("I2").Select()
